$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format while writing, so numeric-looking strings
# like "212.68" or "1.00" are preserved verbatim instead of being coerced
# into numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.960.71'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.642.26'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '212.68'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '23.51'
$ws.Range("D10").Value = '0.0615'
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '1.875.46'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.640.71'
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").Value = '0.571'
$ws.Range("E15").Value = '  +2.09%  '
$ws.Range("D16").Value = '65.56'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '27.944.15'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '233.36'
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").Value = '0.0₃0723'
$ws.Range("D20").Value = '7.61'
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D22").Value = '10.62'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("D25").Value = '153.29'
$ws.Range("E25").Value = '  +2.77%  '
$ws.Range("D26").Value = '6.90'
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("D27").Value = '15.69'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = '0.0485'
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = '3.43'
$ws.Range("E32").Value = '  +4.90%  '
$ws.Range("D33").Value = '3.11'
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").Value = '1.408.51'
$ws.Range("E34").Value = '  -3.70%  '
$ws.Range("E35").Value = '  +2.15%  '
$ws.Range("D36").Value = '2.36'
$ws.Range("E36").Value = '  +1.88%  '
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("D38").Value = '0.566'
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").Value = '0.882'
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").Value = '0.928'
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '67.28'
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("E44").Value = '  +6.43%  '
$ws.Range("D45").Value = '5.52'
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").Value = '1.784.09'
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").Value = '88.10'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").Value = '7.62'
$ws.Range("E51").Value = '  -0.79%  '

# Restore original (unformatted) cell formatting now that the literal text
# values are locked in, so no stray number-format style lingers on the cells.
$ws.Range("D2:D51").ClearFormats()
